$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 8
$ws.Range("H8").Value = 242.18182
$ws.Range("I8").Value = 66.40000000000001
$ws.Range("K8").Value = 199.2
$ws.Range("M8").Value = -60.20000000000002

# ALC row 9
$ws.Range("H9").Value = 19230814
$ws.Range("I9").Value = 20833374
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 20833374
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -20833205
$ws.Range("N9").Value = -438

# ALC row 32
$ws.Range("H32").Value = 13184.5
$ws.Range("I32").Value = 10554.5
$ws.Range("J32").Value = 14499.5
$ws.Range("K32").Value = 10554.5
$ws.Range("L32").Value = 14499.5
$ws.Range("M32").Value = -10228.5
$ws.Range("N32").Value = -15151.5

# ALC row 51
$ws.Range("H51").Value = 8730.5
$ws.Range("J51").Value = 10000
$ws.Range("L51").Value = 10000
$ws.Range("N51").Value = -10968

# ALC row 55
$ws.Range("H55").Value = 612.5
$ws.Range("I55").Value = 700
$ws.Range("J55").Value = 525
$ws.Range("K55").Value = 700
$ws.Range("L55").Value = 525
$ws.Range("M55").Value = -486
$ws.Range("N55").Value = -953

# ALC row 70
$ws.Range("H70").Value = 603328.25
$ws.Range("I70").Value = 929907.4
$ws.Range("J70").Value = 4599.8335
$ws.Range("K70").Value = 2789722.2
$ws.Range("L70").Value = 13799.5005
$ws.Range("M70").Value = -2789452.2
$ws.Range("N70").Value = -14339.5005

# ALC row 73
$ws.Range("H73").Value = 603328.25
$ws.Range("I73").Value = 929907.4
$ws.Range("J73").Value = 4599.8335
$ws.Range("K73").Value = 2789722.2
$ws.Range("L73").Value = 13799.5005
$ws.Range("M73").Value = -2788786.2
$ws.Range("N73").Value = -15671.5005

# ALC row 92
$ws.Range("H92").Value = 37524.668
$ws.Range("I92").Value = 465.25
$ws.Range("K92").Value = 465.25
$ws.Range("M92").Value = 782.75

# ALC row 137
$ws.Range("H137").Value = 4507.8
$ws.Range("I137").Value = 3840.7
$ws.Range("K137").Value = 11522.1
$ws.Range("M137").Value = -8972.099999999999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 3479.1594
$ws.Range("I32").Value = 3309.7354
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 3309.7354
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -3022.7354
$ws.Range("N32").Value = -15574

# ARM row 61
$ws.Range("H61").Value = 2934.1785
$ws.Range("I61").Value = 1588.7368
$ws.Range("K61").Value = 1588.7368
$ws.Range("M61").Value = -1376.7368

# ARM row 97
$ws.Range("H97").Value = 1523.2632
$ws.Range("I97").Value = 1601.2
$ws.Range("J97").Value = 1231
$ws.Range("K97").Value = 1601.2
$ws.Range("L97").Value = 1231
$ws.Range("M97").Value = -1105.2
$ws.Range("N97").Value = -2223

# ARM row 122
$ws.Range("H122").Value = 2047.2354
$ws.Range("I122").Value = 1542.6177
$ws.Range("K122").Value = 4627.8531
$ws.Range("M122").Value = -2177.8531

# ARM row 132
$ws.Range("H132").Value = 6216.5454
$ws.Range("I132").Value = 4536.6665
$ws.Range("J132").Value = 6846.5
$ws.Range("K132").Value = 13609.9995
$ws.Range("L132").Value = 20539.5
$ws.Range("M132").Value = -11079.9995
$ws.Range("N132").Value = -25599.5

# ARM row 136
$ws.Range("H136").Value = 2934.1785
$ws.Range("I136").Value = 1588.7368
$ws.Range("K136").Value = 4766.2104
$ws.Range("M136").Value = -2216.2104

$ws = $wb.Worksheets.Item("BSM")
# BSM row 82
$ws.Range("H82").Value = 13739.6
$ws.Range("I82").Value = 3550
$ws.Range("K82").Value = 3550
$ws.Range("M82").Value = -3167

# BSM row 85
$ws.Range("H85").Value = 13739.6
$ws.Range("I85").Value = 3550
$ws.Range("K85").Value = 3550
$ws.Range("M85").Value = -2224

# BSM row 94
$ws.Range("I94").Value = 12501299
$ws.Range("J94").Value = 1976.1818
$ws.Range("K94").Value = 12501299
$ws.Range("L94").Value = 1976.1818
$ws.Range("M94").Value = -12500848
$ws.Range("N94").Value = -2878.1818

# BSM row 97
$ws.Range("H97").Value = 2799.5
$ws.Range("I97").Value = 2799.5
$ws.Range("K97").Value = 2799.5
$ws.Range("M97").Value = -1808.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 10
$ws.Range("H10").Value = 253
$ws.Range("I10").Value = 253
$ws.Range("K10").Value = 253
$ws.Range("M10").Value = -114

# CRP row 31
$ws.Range("H31").Value = 2393.3906
$ws.Range("J31").Value = 3036.4055
$ws.Range("L31").Value = 3036.4055
$ws.Range("N31").Value = -3626.4055

# CRP row 34
$ws.Range("H34").Value = 2393.3906
$ws.Range("J34").Value = 3036.4055
$ws.Range("L34").Value = 3036.4055
$ws.Range("N34").Value = -3440.4055

# CRP row 59
$ws.Range("H59").Value = 15711.75
$ws.Range("I59").Value = 604
$ws.Range("J59").Value = 17870
$ws.Range("K59").Value = 604
$ws.Range("L59").Value = 17870
$ws.Range("M59").Value = 541
$ws.Range("N59").Value = -20160

$ws = $wb.Worksheets.Item("CUL")
# CUL row 49
$ws.Range("H49").Value = 1666.6666
$ws.Range("I49").Value = 500
$ws.Range("J49").Value = 4000
$ws.Range("K49").Value = 1500
$ws.Range("L49").Value = 12000
$ws.Range("M49").Value = -1344
$ws.Range("N49").Value = -12312

# CUL row 113
$ws.Range("H113").Value = 974.5599999999999
$ws.Range("I113").Value = 496.25
$ws.Range("K113").Value = 1488.75
$ws.Range("M113").Value = 681.25

# CUL row 128
$ws.Range("H128").Value = 250000
$ws.Range("I128").Value = 250000
$ws.Range("K128").Value = 750000
$ws.Range("M128").Value = -745020

$ws = $wb.Worksheets.Item("GSM")
# GSM row 4
$ws.Range("H4").Value = 9999
$ws.Range("I4").Value = 9999
$ws.Range("K4").Value = 9999
$ws.Range("M4").Value = -9887

# GSM row 18
$ws.Range("H18").Value = 7000
$ws.Range("J18").Value = 6250
$ws.Range("L18").Value = 6250
$ws.Range("N18").Value = -6836

# GSM row 109
$ws.Range("H109").Value = 69999
$ws.Range("I109").Value = 69999
$ws.Range("K109").Value = 69999
$ws.Range("M109").Value = -68959

# GSM row 113
$ws.Range("H113").Value = 3471.2942
$ws.Range("I113").Value = 2686.111
$ws.Range("K113").Value = 2686.111
$ws.Range("M113").Value = -516.1109999999999

# GSM row 132
$ws.Range("H132").Value = 3578.0425
$ws.Range("I132").Value = 3580.5757
$ws.Range("K132").Value = 10741.7271
$ws.Range("M132").Value = -8211.7271

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 20003382
$ws.Range("I7").Value = 38464240
$ws.Range("J7").Value = 4118.3335
$ws.Range("K7").Value = 38464240
$ws.Range("L7").Value = 4118.3335
$ws.Range("M7").Value = -38464128
$ws.Range("N7").Value = -4342.3335

# LTW row 68
$ws.Range("H68").Value = 5923.4
$ws.Range("I68").Value = 3930.7778
$ws.Range("K68").Value = 3930.7778
$ws.Range("M68").Value = -3181.7778

# LTW row 71
$ws.Range("H71").Value = 5923.4
$ws.Range("I71").Value = 3930.7778
$ws.Range("K71").Value = 19653.889
$ws.Range("M71").Value = -15909.889

# LTW row 126
$ws.Range("H126").Value = 20003382
$ws.Range("I126").Value = 38464240
$ws.Range("J126").Value = 4118.3335
$ws.Range("K126").Value = 115392720
$ws.Range("L126").Value = 12355.0005
$ws.Range("M126").Value = -115390250
$ws.Range("N126").Value = -17295.0005

# LTW row 132
$ws.Range("H132").Value = 31256202
$ws.Range("I132").Value = 45457956
$ws.Range("K132").Value = 136373868
$ws.Range("M132").Value = -136371338

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Range("H122").Value = 1898.1538
$ws.Range("I122").Value = 1781.08
$ws.Range("K122").Value = 5343.24
$ws.Range("M122").Value = -2893.24

# WVR row 132
$ws.Range("H132").Value = 6522.1665
$ws.Range("I132").Value = 6695.4546
$ws.Range("J132").Value = 6045.625
$ws.Range("K132").Value = 20086.3638
$ws.Range("L132").Value = 18136.875
$ws.Range("M132").Value = -17556.3638
$ws.Range("N132").Value = -23196.875
